$d = $word.ActiveDocument

# Locate the "8/6/2019" date text that follows "Date:" so the edit is
# anchored on content rather than a hard-coded character offset.
$findRng = $d.Content
$found = $findRng.Find.Execute("8/6/2019", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the date text '8/6/2019' to edit."
}
$dateStart = $findRng.Start   # character offset of the leading "8"

# Word marks the point of the most recent edit with a hidden "_GoBack"
# bookmark, splitting the run(s) around it. To reproduce that run layout
# here we first drop a throwaway bookmark right before the "8" so the
# leading space and the date digits land in separate runs, make the
# actual text edit (8 -> 9), then drop the real "_GoBack" bookmark right
# after the newly typed "9" (this also removes any pre-existing "_GoBack"
# bookmark elsewhere in the document, matching Word's singleton
# behaviour), and finally remove the throwaway helper bookmark.
$splitRng = $d.Range($dateStart, $dateStart)
$d.Bookmarks.Add("zzTempSplit", $splitRng) | Out-Null

$digitRng = $d.Range($dateStart, $dateStart + 1)
$digitRng.Text = "9"

$goBackRng = $d.Range($dateStart + 1, $dateStart + 1)
$d.Bookmarks.Add("_GoBack", $goBackRng) | Out-Null

$d.Bookmarks("zzTempSplit").Delete()
